# Update the Maltaspor roster sheet: names, positions and teams are
# reshuffled across rows 2-16 (rows 12, 17, 18 stay unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("Derrick White", "PG,SG", "Boston Celtics")
    3  = @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
    4  = @("Cam Thomas", "SG,SF", "Brooklyn Nets")
    5  = @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans")
    6  = @("Brandon Miller", "SG,SF", "Charlotte Hornets")
    7  = @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
    8  = @("Yves Missi", "C", "New Orleans Pelicans")
    9  = @("Anthony Davis", "PF,C", "Los Angeles Lakers")
    10 = @("Dalton Knecht", "SG", "Los Angeles Lakers")
    11 = @("Gradey Dick", "SG,SF", "Toronto Raptors")
    13 = @("Damian Lillard", "PG", "Milwaukee Bucks")
    14 = @("Bam Adebayo", "C", "Miami Heat")
    15 = @("Jared McCain", "PG,SG", "Philadelphia 76ers")
    16 = @("Cade Cunningham", "PG,SG", "Detroit Pistons")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
